$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: rename columns from Spanish labels to short codes ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Normalize capitalization of Spanish connector words (de/del/el/la/las/los/y) ---
# --- inside state/municipality names to title case, and fix one typo ---
$ws.Range("B5").Value = 'Pabellón De Arteaga'
$ws.Range("B6").Value = 'Rincón De Romos'
$ws.Range("B7").Value = 'San José De Gracia'
$ws.Range("B11").Value = 'Playas De Rosarito'
$ws.Range("B26").Value = 'Amatenango De La Frontera'
$ws.Range("B42").Value = 'Marqués De Comillas'
$ws.Range("B80").Value = 'San Juan De Sabinas'
$ws.Range("B91").Value = 'Villa De Álvarez'
$ws.Range("A93").Value = 'Ciudad De México'
$ws.Range("B97").Value = 'Cuajimalpa De Morelos'
$ws.Range("B116").Value = 'Nombre De Dios'
$ws.Range("B123").Value = 'San Pedro Del Gallo'
$ws.Range("A130").Value = 'Estado De México'
$ws.Range("B130").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B132").Value = 'Almoloya De Alquisiras'
$ws.Range("B135").Value = 'Atizapán De Zaragoza'
$ws.Range("B144").Value = 'Coacalco De Berriozábal'
$ws.Range("B150").Value = 'Ecatepec De Morelos'
$ws.Range("B152").Value = 'Ixtapan De La Sal'
$ws.Range("B153").Value = 'Ixtapan Del Oro'
$ws.Range("B160").Value = 'Naucalpan De Juárez'
$ws.Range("B165").Value = 'San Felipe Del Progreso'
$ws.Range("B174").Value = 'Tlalnepantla De Baz'
$ws.Range("B178").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B179").Value = 'Villa De Allende'
$ws.Range("B188").Value = 'Apaseo El Alto'
$ws.Range("B202").Value = 'San Francisco Del Rincón'
$ws.Range("B204").Value = 'Silao De La Victoria'
$ws.Range("B207").Value = 'Valle De Santiago'
$ws.Range("B209").Value = 'Acapulco De Juárez'
$ws.Range("B211").Value = 'Alcozauca De Guerrero'
$ws.Range("B214").Value = 'Atenango Del Río'
$ws.Range("B215").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B217").Value = 'Atoyac De Álvarez'
$ws.Range("B218").Value = 'Ayutla De Los Libres'
$ws.Range("B221").Value = 'Buenavista De Cuéllar'
$ws.Range("B222").Value = 'Chilapa De Álvarez'
$ws.Range("B223").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B224").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B227").Value = 'Coyuca De Benítez'
$ws.Range("B228").Value = 'Coyuca De Catalán'
$ws.Range("B232").Value = 'Cuetzala Del Progreso'
$ws.Range("B233").Value = 'Cutzamala De Pinzón'
$ws.Range("B239").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B240").Value = 'Iguala De La Independencia'
$ws.Range("B242").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B243").Value = 'Zihuatanejo De Azueta'
$ws.Range("B245").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B248").Value = 'Mártir De Cuilapan'
$ws.Range("B259").Value = 'Taxco De Alarcón'
$ws.Range("B261").Value = 'Técpan De Galeana'
$ws.Range("B263").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B265").Value = 'Tixtla De Guerrero'
$ws.Range("B267").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B268").Value = 'Tlapa De Comonfort'
$ws.Range("B281").Value = 'Atotonilco De Tula'
$ws.Range("B282").Value = 'Atotonilco El Grande'
$ws.Range("B287").Value = 'Cuautepec De Hinojosa'
$ws.Range("B289").Value = 'Huasca De Ocampo'
$ws.Range("B290").Value = 'Huejutla De Reyes'
$ws.Range("B293").Value = 'Jacala De Ledezma'
$ws.Range("B297").Value = 'Mineral Del Chico'
$ws.Range("B298").Value = 'Mixquiahuala De Juárez'
$ws.Range("B299").Value = 'Molango De Escamilla'
$ws.Range("B300").Value = 'Pachuca De Soto'
$ws.Range("B301").Value = 'Progreso De Obregón'
$ws.Range("B305").Value = 'Tenango De Doria'
$ws.Range("B307").Value = 'Tezontepec De Aldama'
$ws.Range("B310").Value = 'Tulancingo De Bravo'
$ws.Range("B311").Value = 'Zacualtipán De Ángeles'
$ws.Range("B314").Value = 'Acatlán De Juárez'
$ws.Range("B315").Value = 'Ahualulco De Mercado'
$ws.Range("B320").Value = 'Atemajac De Brizuela'
$ws.Range("B323").Value = 'Atotonilco El Alto'
$ws.Range("B325").Value = 'Autlán De Navarro'
$ws.Range("B334").Value = 'Cuautitlán De García Barragán'
$ws.Range("B340").Value = 'Encarnación De Díaz'
$ws.Range("B348").Value = 'Jilotlán De Los Dolores'
$ws.Range("B354").Value = 'La Manzanilla De La Paz'
$ws.Range("B355").Value = 'Lagos De Moreno'
$ws.Range("B361").Value = 'Ojuelos De Jalisco'
$ws.Range("B365").Value = 'San Diego De Alejandría'
$ws.Range("B367").Value = 'San Juan De Los Lagos'
$ws.Range("B368").Value = 'San Martín De Bolaños'
$ws.Range("B370").Value = 'San Miguel El Alto'
$ws.Range("B371").Value = 'San Sebastián Del Oeste'
$ws.Range("B372").Value = 'Santa María De Los Ángeles'
$ws.Range("B373").Value = 'Santa María Del Oro'
$ws.Range("B375").Value = 'Talpa De Allende'
$ws.Range("B376").Value = 'Tamazula De Gordiano'
$ws.Range("B382").Value = 'Teocuitatlán De Corona'
$ws.Range("B383").Value = 'Tepatitlán De Morelos'
$ws.Range("B385").Value = 'Tizapán El Alto'
$ws.Range("B386").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B392").Value = 'Unión De San Antonio'
$ws.Range("B393").Value = 'Unión De Tula'
$ws.Range("B394").Value = 'Valle De Guadalupe'
$ws.Range("B397").Value = 'Zacoalco De Torres'
$ws.Range("B399").Value = 'Zapotitlán De Vadillo'
$ws.Range("B400").Value = 'Zapotlán El Grande'
$ws.Range("B420").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B422").Value = 'Cojumatlán De Régules'
$ws.Range("B476").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B498").Value = 'Coatlán Del Río'
$ws.Range("B506").Value = 'Puente De Ixtla'
$ws.Range("B511").Value = 'Tetela Del Volcán'
$ws.Range("B512").Value = 'Tlaltizapán De Zapata'
$ws.Range("B522").Value = 'Ixtlán Del Río'
$ws.Range("B528").Value = 'Santa María Del Oro'
$ws.Range("B540").Value = 'Montemorelos'
$ws.Range("B542").Value = 'San Nicolás De Los Garza'
$ws.Range("B545").Value = 'Ayoquezco De Aldama'
$ws.Range("B549").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B551").Value = 'Coicoyán De Las Flores'
$ws.Range("B552").Value = 'Constancia Del Rosario'
$ws.Range("B554").Value = 'Cuilápam De Guerrero'
$ws.Range("B555").Value = 'Fresnillo De Trujano'
$ws.Range("B556").Value = 'Guadalupe De Ramírez'
$ws.Range("B557").Value = 'Guevea De Humboldt'
$ws.Range("B558").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B559").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B560").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B561").Value = 'Huautla De Jiménez'
$ws.Range("B563").Value = 'Ixtlán De Juárez'
$ws.Range("B564").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B567").Value = 'Mariscala De Juárez'
$ws.Range("B568").Value = 'Mártires De Tacubaya'
$ws.Range("B571").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B572").Value = 'Mixistlán De La Reforma'
$ws.Range("B574").Value = 'Oaxaca De Juárez'
$ws.Range("B575").Value = 'Ocotlán De Morelos'
$ws.Range("B576").Value = 'Putla Villa De Guerrero'
$ws.Range("B592").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B599").Value = 'San José Del Progreso'
$ws.Range("B606").Value = 'San Juan Del Estado'
$ws.Range("B632").Value = 'San Miguel El Grande'
$ws.Range("B645").Value = 'San Pedro El Alto'
$ws.Range("B661").Value = 'Santa Cruz De Bravo'
$ws.Range("B712").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B713").Value = 'Tataltepec De Valdés'
$ws.Range("B714").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B715").Value = 'Tlacolula De Matamoros'
$ws.Range("B718").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B719").Value = 'Villa De Etla'
$ws.Range("B720").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B721").Value = 'Villa De Zaachila'
$ws.Range("B722").Value = 'Villa Sola De Vega'
$ws.Range("B723").Value = 'Villa Talea De Castro'
$ws.Range("B726").Value = 'Zimatlán De Álvarez'
$ws.Range("B746").Value = 'Cuayuca De Andrade'
$ws.Range("B755").Value = 'Huehuetlán El Chico'
$ws.Range("B756").Value = 'Huehuetlán El Grande'
$ws.Range("B761").Value = 'Izúcar De Matamoros'
$ws.Range("B766").Value = 'Los Reyes De Juárez'
$ws.Range("B773").Value = 'Palmar De Bravo'
$ws.Range("B782").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B784").Value = 'San Salvador El Seco'
$ws.Range("B785").Value = 'San Salvador El Verde'
$ws.Range("B786").Value = 'Tecali De Herrera'
$ws.Range("B791").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B796").Value = 'Tepexi De Rodríguez'
$ws.Range("B797").Value = 'Tetela De Ocampo'
$ws.Range("B800").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B817").Value = 'Amealco De Bonfil'
$ws.Range("B818").Value = 'Cadereyta De Montes'
$ws.Range("B825").Value = 'Pinal De Amoles'
$ws.Range("B827").Value = 'San Juan Del Río'
$ws.Range("B840").Value = 'Santa María Del Río'
$ws.Range("B843").Value = 'Villa De Ramos'
$ws.Range("B899").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B903").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B904").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B906").Value = 'San Pablo Del Monte'
$ws.Range("B923").Value = 'Amatlán De Los Reyes'
$ws.Range("B929").Value = 'Castillo De Teayo'
$ws.Range("B935").Value = 'Cosamaloapan De Carpio'
$ws.Range("B945").Value = 'Hueyapan De Ocampo'
$ws.Range("B946").Value = 'Ignacio De La Llave'
$ws.Range("B949").Value = 'Ixhuatlán De Madero'
$ws.Range("B950").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B963").Value = 'Martínez De La Torre'
$ws.Range("B972").Value = 'Paso De Ovejas'
$ws.Range("B975").Value = 'Poza Rica De Hidalgo'
$ws.Range("B979").Value = 'Sayula De Alemán'
$ws.Range("B981").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1000").Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range("B1022").Value = 'Jiménez Del Teul'
$ws.Range("B1024").Value = 'Mezquital Del Oro'
$ws.Range("B1027").Value = 'Nochistlán De Mejía'
$ws.Range("B1034").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1035").Value = 'Trinidad García De La Cadena'
$ws.Range("B1037").Value = 'Villa De Cos'

# --- Correct floating point rounding for Lázaro Cárdenas percentage ---
$ws.Range("D442").Value = 0.009310618066561017

# --- Remove trailing footnote/metadata rows (now rows 1045-1049) ---
$ws.Range("A1045:A1049").EntireRow.Delete()

